# Chapter 4 final edits
# Shift the "date" column (F) values forward by 2 (serial date numbers),
# for rows 2 through 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value2 = 44578
$ws.Range("F3").Value2 = 44577
$ws.Range("F4").Value2 = 44576
$ws.Range("F5").Value2 = 44575
$ws.Range("F6").Value2 = 44574
$ws.Range("F7").Value2 = 44573
